{"js": "// ONC-2660: DESIGN: Add matrix information to student view of submission\n//\n// Appends two new rows to the end of the Table of Contents table, matching\n// the existing row style (bold, blue, underlined Trebuchet MS heading in the\n// left/wide cell; bold, right-aligned Trebuchet MS page number in the\n// right/narrow cell):\n//   SV: ONE SUBMISSION (SVOS-1)        -> 33\n//   SV: MULTIPLE SUBMISSIONS (SVMS-1)  -> 35\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newRows = [\n  [\"SV: ONE SUBMISSION (SVOS-1)\", \"33\"],\n  [\"SV: MULTIPLE SUBMISSIONS (SVMS-1)\", \"35\"],\n];\n\nfor (const rowValues of newRows) {\n  table.addRows(\"End\", 1, [rowValues]);\n  await context.sync();\n}\n", "ps1": "# ONC-2660: DESIGN: Add matrix information to student view of submission\n#\n# Appends two new rows to the end of the Table of Contents table, matching\n# the existing row style (bold, blue, underlined Trebuchet MS heading in the\n# left/wide cell; bold, right-aligned Trebuchet MS page number in the\n# right/narrow cell):\n#   SV: ONE SUBMISSION (SVOS-1)        -> 33\n#   SV: MULTIPLE SUBMISSIONS (SVMS-1)  -> 35\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rows = @(\n    @(\"SV: ONE SUBMISSION (SVOS-1)\", \"33\"),\n    @(\"SV: MULTIPLE SUBMISSIONS (SVMS-1)\", \"35\")\n)\n\nforeach ($rowData in $rows) {\n    $newRow = $t.Rows.Add()\n    $headingCell = $newRow.Cells.Item(1)\n    $pageCell = $newRow.Cells.Item(2)\n    $headingCell.Range.Text = $rowData[0]\n    $pageCell.Range.Text = $rowData[1]\n}\n"}
